$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(0.2431903771268651,0.2141971996420864,0.1963765498992132,0.1891102528108206,0.1879034476085053,0.1962785706105024,0.233197688394057,0.305430200716188,0.3583802879072095,0.382439243199542,0.3915452584293746,0.3895843274798949,0.3831884954948919,0.3792702507084584,0.3568073646844141,0.3430194930081996,0.335086440712189,0.3324000090688344,0.34448751176771,0.3850672350651223,0.4115614411131219,0.3974236256450752,0.3438238398256317,0.2859090082249338)
$Cvals = @(0.05551934676972081,0.05339687447753505,0.05208753242827413,0.05155246469899311,0.05146352775679475,0.0520803223342341,0.05478881664731716,0.06004985468139523,0.06388263683051321,0.06561884507701166,0.06627520729463754,0.06613389769985645,0.0656728667459987,0.06539032707128456,0.06376901981606409,0.06277248657086432,0.0621986185456791,0.06200419985756866,0.06287864082479189,0.06580831299217493,0.06771657750996951,0.06669870581696102,0.06283065143638566,0.05863218047436192)
$Evals = @(0.5648735511074534,0.5550854854953968,0.5493443845042734,0.5470725275157164,0.5466993783667249,0.5493134712554379,0.5614428862328893,0.5873592613884284,0.6076986032895775,0.6172335779191513,0.6208847970912643,0.620096640558657,0.6175331541448372,0.6159682202974466,0.6070811470916766,0.6017015061817546,0.5986338753872644,0.5975997990841648,0.6022714261866469,0.6182850129802659,0.6289870212859938,0.6232535761511855,0.6020136868605164,0.58012012743648)
$Fvals = @(2.044944038988177,2.033183166965259,2.027051235759373,2.024826485584526,2.024473624618977,2.027020122137841,2.04066288886726,2.076056842888121,2.107330967253333,2.122703716554952,2.12868972423928,2.127393207926616,2.123192889195806,2.120641515716173,2.106349372761727,2.097875043548868,2.09310870655392,2.091513437416424,2.098765986707221,2.124422156117774,2.142149781293114,2.132600426181057,2.098362862148406,2.065556531768394)
$Gvals = @(0.4283065998189315,0.4312419932998282,0.4333232307011912,0.4342414385632907,0.4343981380734334,0.4333353302530583,0.4292608125124531,0.4234858263913566,0.4205968319777895,0.4195773131610423,0.4192336825448848,0.4193058010828139,0.4195481915708612,0.419702191379379,0.4206693953126788,0.4213382671530823,0.421750719375396,0.4218951302095562,0.4212641937047863,0.4194758433366559,0.4185544484137793,0.4190235572862022,0.4212975953811195,0.4248106067133151)
$Hvals = @(0.5961500072058001,0.6006553913787229,0.6036540412581672,0.6049344814160662,0.6051506300807148,0.6036710729349295,0.5976552719521777,0.5876997083025941,0.5815052614804799,0.5789298730355412,0.5779894699062709,0.5781904534783493,0.5788518074180757,0.5792614421866915,0.581678446239593,0.5832232875799619,0.5841346661556699,0.584447164712607,0.5830564743233921,0.5786566064100711,0.5759841021945178,0.5773918975106582,0.5831318182382432,0.5901960938229251)
$Jvals = @(0.03999789559807532,0.04009112975882978,0.04015490912184116,0.04018254748266159,0.04018723647642375,0.04015527518395778,0.04002868971353912,0.03983208096659752,0.03971882258638004,0.03967401170004869,0.03965800308625944,0.039661408193016,0.03967267544297925,0.039679701877537,0.03972188572016222,0.0397494796037936,0.03976598282094557,0.03977167920804803,0.03974647681774712,0.03966933995387656,0.03962452171001551,0.03964793165214431,0.03974783238659896,0.03987977120494435)
$Kvals = @(0.2223516965622281,0.1939576610440668,0.1764592190571506,0.169312717384912,0.1681251084514486,0.1763629019932154,0.2125750567061004,0.283059719070053,0.3345061854599862,0.3578336274262028,0.3666558267202333,0.3647563230002788,0.3585596663120327,0.3547625362969029,0.3329801144002147,0.3195975342827353,0.3118931217790077,0.3092833346388204,0.3210228718496069,0.3603800878364041,0.3860356454392502,0.3723490541737817,0.3203785096814613,0.2640499746988212)
$Mvals = @(0.4018644365232333,0.3811933185100642,0.3686580187189108,0.3635894776062401,0.3627502551963104,0.3685895015396099,0.3947046262561997,0.4471525396398803,0.4864327252334775,0.5044632134837315,0.5113139463295369,0.5098375009229059,0.5050263684806069,0.502082395093808,0.4852576211085591,0.4749773931279222,0.4690797178718213,0.4670854923351797,0.4760701656844333,0.5064388933068997,0.5264204171781302,0.5157437510867666,0.4755760840530314,0.4328323565412759)
$Nvals = @(1.413611318524705,1.427288730535384,1.436149594181771,1.439876902238666,1.440502852452287,1.436199390470186,1.418231244861438,1.386666291691107,1.365709556347337,1.35666014457184,1.35330292302017,1.354022865366939,1.356382549104822,1.357836987047595,1.366310694893798,1.371632986856383,1.37473977587868,1.375799503867892,1.371061705746971,1.355687563944276,1.346045327876475,1.351154449962266,1.37131983574438,1.39481287650435)
$Ovals = @(1.999951363741616,2.015394959385588,2.025947734165626,2.030517221244651,2.031292236815375,2.026008270107781,2.005054163807046,1.972456707558763,1.953687339506587,1.946274009960604,1.943628581507198,1.944191122901032,1.946053125072396,1.947214733456605,1.954194461071893,1.958764476034446,1.961498909108201,1.962442924184515,1.958267031293602,1.94550181707497,1.938102363753416,1.941965250759665,1.958491592548,1.98036567762685)

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Range("B" + $r).Value = $Bvals[$i]
    $ws.Range("C" + $r).Value = $Cvals[$i]
    $ws.Range("E" + $r).Value = $Evals[$i]
    $ws.Range("F" + $r).Value = $Fvals[$i]
    $ws.Range("G" + $r).Value = $Gvals[$i]
    $ws.Range("H" + $r).Value = $Hvals[$i]
    $ws.Range("J" + $r).Value = $Jvals[$i]
    $ws.Range("K" + $r).Value = $Kvals[$i]
    $ws.Range("M" + $r).Value = $Mvals[$i]
    $ws.Range("N" + $r).Value = $Nvals[$i]
    $ws.Range("O" + $r).Value = $Ovals[$i]
}
